# Swap the "B" and "C" quarter rows within every year group (columns A:E),
# then delete the now-redundant F and G columns (which duplicated B/E with
# rounding-diff / cumulative-to-period conversions no longer needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
  ,@(3,4)
  ,@(7,8)
  ,@(11,12)
  ,@(15,16)
  ,@(19,20)
  ,@(23,24)
  ,@(27,28)
  ,@(31,32)
  ,@(35,36)
  ,@(39,40)
  ,@(43,44)
  ,@(47,48)
  ,@(51,52)
  ,@(55,56)
  ,@(59,60)
  ,@(63,64)
)

foreach ($pair in $rowPairs) {
  $r1 = $pair[0]
  $r2 = $pair[1]
  $rangeA = $ws.Range("A" + $r1 + ":E" + $r1)
  $rangeB = $ws.Range("A" + $r2 + ":E" + $r2)
  $valA = $rangeA.Value()
  $valB = $rangeB.Value()
  $rangeA.Value = $valB
  $rangeB.Value = $valA
}

$ws.Range("F1:G65").EntireColumn.Delete()
